$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster was refreshed: two players (Derrick Jones Jr., Jalen Smith) were
# dropped, one player (Nicolas Claxton) was added, and the remaining players
# were re-sorted into a new row order (each player keeps its own
# position/team alongside it). Rewrite the whole A2:C18 block row by row,
# then trim the now-unused 19th row.

$rows = @(
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Kel'el Ware", "PF,C", "Miami Heat"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# Row 19 (formerly Alperen Sengün, now re-seated at row 11) is no longer
# needed, shrinking the sheet from 19 to 18 rows.
$ws.Rows.Item(19).Delete()
